$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: matricula changes, keep it stored as text (leading/ID-like value)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "123456"

# Row 3: new student data (matricula, nome, and grades)
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "000111"
$ws.Range("B3").Value = "thiago"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 2.5
